$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (one month later: 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update prices in column D (rows 14-21)
$ws.Range("D14").Value = 83.175
$ws.Range("D15").Value = 108.235
$ws.Range("D16").Value = 129.622
$ws.Range("D17").Value = 206.747
$ws.Range("D18").Value = 257.947
$ws.Range("D19").Value = 322.11
$ws.Range("D20").Value = 360.78
$ws.Range("D21").Value = 399.666
